$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$cellValues = @{
    "B2" = 1.02
    "C2" = 1.039806583996589
    "D2" = 1.04316315749019
    "E2" = 1.043461600854908
    "F2" = 1.038549137735114
    "I2" = 1.041490219652241
    "J2" = 1.044896605493829
    "K2" = 1.045937798595537
    "L2" = 1.046235400509971
    "M2" = 1.041336854701452
    "N2" = 1.018621321845196
    "B3" = 1.02
    "C3" = 1.041124566061259
    "D3" = 1.044207470626641
    "E3" = 1.044728947576123
    "F3" = 1.040497210900389
    "I3" = 1.041956922537177
    "J3" = 1.045858103701232
    "K3" = 1.046792480451635
    "L3" = 1.04731259435626
    "M3" = 1.0430919612672
    "N3" = 1.018965337162738
    "B4" = 1.02
    "C4" = 1.041975569059754
    "D4" = 1.044881597741098
    "E4" = 1.04554754924184
    "F4" = 1.04175564747254
    "I4" = 1.042256634544827
    "J4" = 1.046477986730183
    "K4" = 1.047343307354285
    "L4" = 1.048007606549122
    "M4" = 1.044225143340536
    "N4" = 1.019186566634232
    "B5" = 1.02
    "C5" = 1.042332902217841
    "D5" = 1.045164619478563
    "E5" = 1.045891346996114
    "F5" = 1.042284207437298
    "I5" = 1.042382091610713
    "J5" = 1.046738047493262
    "K5" = 1.047574350784626
    "L5" = 1.04829931568133
    "M5" = 1.044700950953768
    "N5" = 1.019279245076896
    "B6" = 1.02
    "C6" = 1.042392875105104
    "D6" = 1.04521211785641
    "E6" = 1.045949052255802
    "F6" = 1.042372926878221
    "I6" = 1.042403124716014
    "J6" = 1.046781681447226
    "K6" = 1.047613113390844
    "L6" = 1.048348267301548
    "M6" = 1.044780807408368
    "N6" = 1.019294787108305
    "B7" = 1.02
    "C7" = 1.041980345437861
    "D7" = 1.04488538098359
    "E7" = 1.04555214442185
    "F7" = 1.041762712010526
    "I7" = 1.042258313033732
    "J7" = 1.046481463782239
    "K7" = 1.047346396621825
    "L7" = 1.048011506236111
    "M7" = 1.044231503375536
    "N7" = 1.019187806288405
    "B8" = 1.02
    "C8" = 1.040252383237523
    "D8" = 1.043516425031384
    "E8" = 1.043890211966405
    "F8" = 1.039207941603422
    "I8" = 1.041648417149662
    "J8" = 1.045222021796532
    "K8" = 1.046227102666123
    "L8" = 1.046599861729788
    "M8" = 1.041930524928049
    "N8" = 1.018737868545882
    "B9" = 1.02
    "C9" = 1.037193238113867
    "D9" = 1.041091571224017
    "E9" = 1.040950237278742
    "F9" = 1.034689362504233
    "I9" = 1.04055613573744
    "J9" = 1.042985084988655
    "K9" = 1.044237628612108
    "L9" = 1.044096750858275
    "M9" = 1.037856221685946
    "N9" = 1.017934431530194
    "B10" = 1.02
    "C10" = 1.03514377420464
    "D10" = 1.039466230873729
    "E10" = 1.038982165626064
    "F10" = 1.031664754266896
    "I10" = 1.039815945484323
    "J10" = 1.041481583496918
    "K10" = 1.042899490555881
    "L10" = 1.04241714300206
    "M10" = 1.035125914009289
    "N10" = 1.017391566399371
    "B11" = 1.02
    "C11" = 1.034253852941645
    "D11" = 1.038760289366888
    "E11" = 1.038127960787723
    "F11" = 1.030351951797182
    "I11" = 1.039492544873438
    "J11" = 1.040827575815558
    "K11" = 1.042317190142866
    "L11" = 1.041687193523251
    "M11" = 1.033940120676688
    "N11" = 1.017154755217247
    "B12" = 1.02
    "C12" = 1.033922913699335
    "D12" = 1.038497740782262
    "E12" = 1.037810360255877
    "F12" = 1.02986383012608
    "I12" = 1.039371980997601
    "J12" = 1.040584193403401
    "K12" = 1.04210045954529
    "L12" = 1.041415649658298
    "M12" = 1.033499113762731
    "N12" = 1.017066528136097
    "B13" = 1.02
    "C13" = 1.033993918757466
    "D13" = 1.03855407339993
    "E13" = 1.037878500839715
    "F13" = 1.029968556281745
    "I13" = 1.039397862262178
    "J13" = 1.040636420450717
    "K13" = 1.042146968924063
    "L13" = 1.041473915325088
    "M13" = 1.033593736424504
    "N13" = 1.017085465168893
    "B14" = 1.02
    "C14" = 1.034226505262022
    "D14" = 1.038738593796709
    "E14" = 1.038101714214834
    "F14" = 1.030311613578445
    "I14" = 1.039482587993007
    "J14" = 1.040807467068109
    "K14" = 1.042299284099767
    "L14" = 1.041664755990603
    "M14" = 1.033903678234922
    "N14" = 1.017147467756293
    "B15" = 1.02
    "C15" = 1.034369758486652
    "D15" = 1.038852238925101
    "E15" = 1.038239201920226
    "F15" = 1.030522917227543
    "I15" = 1.039534732093705
    "J15" = 1.040912794056521
    "K15" = 1.042393072249188
    "L15" = 1.04178228492385
    "M15" = 1.034094570164971
    "N15" = 1.017185634427562
    "B16" = 1.02
    "C16" = 1.035202782141205
    "D16" = 1.039513035883626
    "E16" = 1.039038813230327
    "F16" = 1.031751813177523
    "I16" = 1.039837347235076
    "J16" = 1.04152492442743
    "K16" = 1.042938074778898
    "L16" = 1.042465530510977
    "M16" = 1.035204534937774
    "N16" = 1.01740724576626
    "B17" = 1.02
    "C17" = 1.035724643085633
    "D17" = 1.039926954133056
    "E17" = 1.039539842823904
    "F17" = 1.032521816866183
    "I17" = 1.040026392418609
    "J17" = 1.041908094533875
    "K17" = 1.043279165822532
    "L17" = 1.042893392961091
    "M17" = 1.035899823701981
    "N17" = 1.017545787297692
    "B18" = 1.02
    "C18" = 1.036028796107694
    "D18" = 1.040168177812598
    "E18" = 1.039831890847182
    "F18" = 1.032970646529183
    "I18" = 1.04013638036856
    "J18" = 1.042131303960805
    "K18" = 1.043477841195917
    "L18" = 1.043142700942109
    "M18" = 1.036305032293195
    "N18" = 1.017626427832425
    "B19" = 1.02
    "C19" = 1.036132464045861
    "D19" = 1.040250393787892
    "E19" = 1.039931438994707
    "F19" = 1.033123635440953
    "I19" = 1.04017383622379
    "J19" = 1.042207364097851
    "K19" = 1.043545537514209
    "L19" = 1.043227665156609
    "M19" = 1.036443140519236
    "N19" = 1.017653895679919
    "B20" = 1.02
    "C20" = 1.035668677197332
    "D20" = 1.039882566165482
    "E20" = 1.039486107221128
    "F20" = 1.032439233973016
    "I20" = 1.040006138536496
    "J20" = 1.041867013766805
    "K20" = 1.043242598704303
    "L20" = 1.042847514000748
    "M20" = 1.035825261253494
    "N20" = 1.017530940538449
    "B21" = 1.02
    "C21" = 1.03415802497333
    "D21" = 1.038684266311809
    "E21" = 1.038035992082655
    "F21" = 1.030210605386742
    "I21" = 1.039457650512154
    "J21" = 1.040757110713949
    "K21" = 1.042254443224031
    "L21" = 1.041608569496179
    "M21" = 1.033812423428417
    "N21" = 1.017129216870591
    "B22" = 1.02
    "C22" = 1.033205998432541
    "D22" = 1.037928932477895
    "E22" = 1.03712244489702
    "F22" = 1.028806543192029
    "I22" = 1.039110255448279
    "J22" = 1.040056635591902
    "K22" = 1.041630611815696
    "L22" = 1.04082723096476
    "M22" = 1.032543679088463
    "N22" = 1.016875103630445
    "B23" = 1.02
    "C23" = 1.033710899055229
    "D23" = 1.038329532855843
    "E23" = 1.037606907005805
    "F23" = 1.029551138091992
    "I23" = 1.039294658077876
    "J23" = 1.040428222793173
    "K23" = 1.041961559299456
    "L23" = 1.04124165991591
    "M23" = 1.033216572812288
    "N23" = 1.017009960043199
    "B24" = 1.02
    "C24" = 1.035693966512912
    "D24" = 1.039902623819032
    "E24" = 1.039510388625963
    "F24" = 1.03247655055416
    "I24" = 1.040015291255448
    "J24" = 1.041885577283978
    "K24" = 1.043259122667196
    "L24" = 1.042868245521386
    "M24" = 1.035858953868163
    "N24" = 1.017537649670009
    "B25" = 1.02
    "C25" = 1.037985835782106
    "D25" = 1.041719976809376
    "E25" = 1.041711687408803
    "F25" = 1.035859606688001
    "I25" = 1.040840618046808
    "J25" = 1.043565513003633
    "K25" = 1.044754016779552
    "L25" = 1.044745753108023
    "M25" = 1.038911948656513
    "N25" = 1.018143406508395
}

foreach ($cellRef in $cellValues.Keys) {
    $ws.Range($cellRef).Value = $cellValues[$cellRef]
}
